$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure every Price cell we touch keeps its original plain-text storage
# (these look numeric to Excel, so force Text format before writing the value,
# otherwise COM auto-converts "7.30" -> 7.3, "0.0000225" -> 2.25E-05, etc.)
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D11", "D12", "D13", "D14", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D31", "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "61.441.89"
$ws.Range("E2").Value = "  -2.37%  "
$ws.Range("D3").Value = "2.969.21"
$ws.Range("E3").Value = "  -2.51%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "586.33"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").Value = "140.92"
$ws.Range("E6").Value = "  -6.47%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.518"
$ws.Range("E8").Value = "  -2.32%  "
$ws.Range("D9").Value = "2.968.90"
$ws.Range("E9").Value = "  -2.60%  "
$ws.Range("E10").Value = "  -6.30%  "
$ws.Range("D11").Value = "5.74"
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").Value = "  +2.41%  "
$ws.Range("D13").Value = "0.0000225"
$ws.Range("E13").Value = "  -3.27%  "
$ws.Range("D14").Value = "33.88"
$ws.Range("E14").Value = "  -5.59%  "
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").Value = "3.460.27"
$ws.Range("E16").Value = "  -2.41%  "
$ws.Range("D17").Value = "6.98"
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").Value = "61.477.70"
$ws.Range("E18").Value = "  -2.15%  "
$ws.Range("D19").Value = "2.970.01"
$ws.Range("E19").Value = "  -2.41%  "
$ws.Range("D20").Value = "448.12"
$ws.Range("E20").Value = "  -6.50%  "
$ws.Range("D21").Value = "13.82"
$ws.Range("E21").Value = "  -2.91%  "
$ws.Range("D22").Value = "0.681"
$ws.Range("E22").Value = "  -2.93%  "
$ws.Range("D23").Value = "7.30"
$ws.Range("E23").Value = "  -2.38%  "
$ws.Range("D24").Value = "81.07"
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("D25").Value = "12.06"
$ws.Range("E25").Value = "  -4.18%  "
$ws.Range("D26").Value = "2.14"
$ws.Range("E26").Value = "  -9.78%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "9.79"
$ws.Range("E28").Value = "  -6.96%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("D31").Value = "6.82"
$ws.Range("E31").Value = "  -7.00%  "
$ws.Range("D32").Value = "2.05"
$ws.Range("E32").Value = "  -6.01%  "
$ws.Range("D33").Value = "27.01"
$ws.Range("E33").Value = "  -2.22%  "
$ws.Range("D34").Value = "0.106"
$ws.Range("E34").Value = "  -3.12%  "
$ws.Range("E35").Value = "  -4.48%  "
$ws.Range("D36").Value = "0.0₃0770"
$ws.Range("E36").Value = "  -4.54%  "
$ws.Range("D37").Value = "5.67"
$ws.Range("E37").Value = "  -3.30%  "
$ws.Range("D38").Value = "2.06"
$ws.Range("E38").Value = "  -5.28%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "50.06"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").Value = "9.11"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").Value = "0.118"
$ws.Range("E41").Value = "  +3.37%  "
$ws.Range("D42").Value = "2.77"
$ws.Range("E42").Value = "  -10.95%  "
$ws.Range("D43").Value = "386.71"
$ws.Range("E43").Value = "  -8.93%  "
$ws.Range("D44").Value = "0.0352"
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("D45").Value = "2.719.19"
$ws.Range("E45").Value = "  -3.85%  "
$ws.Range("E46").Value = "  -7.81%  "
$ws.Range("D47").Value = "36.92"
$ws.Range("E47").Value = "  -1.88%  "
$ws.Range("D48").Value = "129.79"
$ws.Range("E48").Value = "  +2.56%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").Value = "0.107"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").Value = "2.15"
$ws.Range("E51").Value = "  -0.92%  "
